$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 290-291; this shifts the existing rows
# 290..378 down to 292..380 (all their data, including per-row dates,
# quality, volume, prices, etc. move along with them - no further edits
# needed for those rows).
$ws.Range("A290:A291").EntireRow.Insert()

# --- Fill in the first new row (290) ---
$ws.Cells.Item(290, 1).Value = 5
$ws.Cells.Item(290, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(290, 3).Value = "Maule"
$ws.Cells.Item(290, 4).Value = 44809
$ws.Cells.Item(290, 5).Value = 7
$ws.Cells.Item(290, 6).Value = 100114014
$ws.Cells.Item(290, 7).Value = "Betarraga"
$ws.Cells.Item(290, 8).Value = "Sin especificar"
$ws.Cells.Item(290, 9).Value = "Primera"
$ws.Cells.Item(290, 10).Value = 3000
$ws.Cells.Item(290, 11).Value = 1000
$ws.Cells.Item(290, 12).Value = 1000
$ws.Cells.Item(290, 13).Value = 1000
$ws.Cells.Item(290, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(290, 15).Value = "Región del Maule"
$ws.Cells.Item(290, 16).Value = 200
$ws.Cells.Item(290, 17).Value = 5
$ws.Cells.Item(290, 18).Value = "Hortaliza"

# --- Fill in the second new row (291) ---
$ws.Cells.Item(291, 1).Value = 5
$ws.Cells.Item(291, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(291, 3).Value = "Maule"
$ws.Cells.Item(291, 4).Value = 44809
$ws.Cells.Item(291, 5).Value = 7
$ws.Cells.Item(291, 6).Value = 100114014
$ws.Cells.Item(291, 7).Value = "Betarraga"
$ws.Cells.Item(291, 8).Value = "Sin especificar"
$ws.Cells.Item(291, 9).Value = "Segunda"
$ws.Cells.Item(291, 10).Value = 2000
$ws.Cells.Item(291, 11).Value = 800
$ws.Cells.Item(291, 12).Value = 800
$ws.Cells.Item(291, 13).Value = 800
$ws.Cells.Item(291, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(291, 15).Value = "Región del Maule"
$ws.Cells.Item(291, 16).Value = 160
$ws.Cells.Item(291, 17).Value = 5
$ws.Cells.Item(291, 18).Value = "Hortaliza"
